# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps with freshly generated values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 21911181-...md
$wsOverview.Range("G2").Value = "2016-09-03 11:10:22"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-03 11:10:18"
$wsZhCn.Range("K2").Value = "2016-09-03 11:10:35"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-03 11:10:22"
$wsDeDe.Range("K2").Value = "2016-09-03 11:10:42"
